$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 220, shifting existing rows 220-223 down to 221-224
$ws.Rows.Item(220).Insert()

# Fill in the new row 220 with data
$ws.Cells.Item(220, 1).Value = 4
$ws.Cells.Item(220, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(220, 3).Value = 'Los Lagos'
$ws.Cells.Item(220, 4).Value = 44656
$ws.Cells.Item(220, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Cells.Item(220, 5).Value = 10
$ws.Cells.Item(220, 6).Value = 100112032
$ws.Cells.Item(220, 7).Value = 'Zapallo italiano'
$ws.Cells.Item(220, 8).Value = 'Sin especificar'
$ws.Cells.Item(220, 9).Value = 'Primera'
$ws.Cells.Item(220, 10).Value = 200
$ws.Cells.Item(220, 11).Value = 13000
$ws.Cells.Item(220, 12).Value = 13000
$ws.Cells.Item(220, 13).Value = 13000
$ws.Cells.Item(220, 14).Value = '$/caja 50 unidades'
$ws.Cells.Item(220, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(220, 16).Value = 260
$ws.Cells.Item(220, 17).Value = 50
$ws.Cells.Item(220, 18).Value = 'Hortaliza'
